$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.447.72"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").Value = "1.573.48"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3733"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.92%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3393"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07553"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.136"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.68%  "

$ws.Range("E14").Value = "  -0.40%  "

$ws.Range("E15").Value = "  -0.08%  "

$ws.Range("D16").Value = "1.569.03"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001124"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06737"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.282"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.52%  "

$ws.Range("D24").Value = "22.443.62"
$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.336"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.607"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.98%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.016"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.52%  "

$ws.Range("D31").Value = "1.747.28"
$ws.Range("E31").Value = "  +0.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.050"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.133"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.981"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.787"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08389"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.379"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02463"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2289"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06506"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.01%  "

$ws.Range("E41").Value = "  -0.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6223"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.809"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5807"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.52%  "

$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("E50").Value = "  -6.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07316"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "
